$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.979788333333333
$ws.Range("H2").Value = 14.939365
$ws.Range("I2").Value = 0.129176854764059
$ws.Range("J2").Value = 0.129176854764059
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 118.0470123333333
$ws.Range("N2").Value = 354.141037
$ws.Range("O2").Value = 0.4657216250363638
$ws.Range("P2").Value = 0.4657216250363638
$ws.Range("Q2").Value = 587.8491348023895
$ws.Range("R2").Value = 5290.642213221505
$ws.Range("S2").Value = 0.06016045471780389
$ws.Range("T2").Value = 0.06016045471780391
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.979788333333333
$ws.Range("H3").Value = 14.939365
$ws.Range("I3").Value = 0.129176854764059
$ws.Range("J3").Value = 0.129176854764059
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2342114953037475
$ws.Range("P3").Value = 0.2342114953037476
$ws.Range("Q3").Value = 295.6294435851711
$ws.Range("R3").Value = 2660.66499226654
$ws.Range("S3").Value = 0.03025470431292528
$ws.Range("T3").Value = 0.03025470431292529
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.979788333333333
$ws.Range("H4").Value = 14.939365
$ws.Range("I4").Value = 0.129176854764059
$ws.Range("J4").Value = 0.129176854764059
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 56.84506433333333
$ws.Range("N4").Value = 170.535193
$ws.Range("O4").Value = 0.2242663767030476
$ws.Range("P4").Value = 0.2242663767030477
$ws.Range("Q4").Value = 283.0763881747161
$ws.Range("R4").Value = 2547.687493572445
$ws.Range("S4").Value = 0.02897002517183132
$ws.Range("T4").Value = 0.02897002517183133
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.979788333333333
$ws.Range("H5").Value = 14.939365
$ws.Range("I5").Value = 0.129176854764059
$ws.Range("J5").Value = 0.129176854764059
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.21324333333333
$ws.Range("N5").Value = 57.63973
$ws.Range("O5").Value = 0.07580050295684103
$ws.Range("P5").Value = 0.07580050295684104
$ws.Range("Q5").Value = 95.67788499682779
$ws.Range("R5").Value = 861.10096497145
$ws.Range("S5").Value = 0.009791670561498475
$ws.Range("T5").Value = 0.00979167056149848
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.019504
$ws.Range("H6").Value = 54.058512
$ws.Range("I6").Value = 0.467430078412646
$ws.Range("J6").Value = 0.4674300784126461
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.0470123333333
$ws.Range("N6").Value = 354.141037
$ws.Range("O6").Value = 0.4657216250363638
$ws.Range("P6").Value = 0.4657216250363638
$ws.Range("Q6").Value = 2127.148610928549
$ws.Range("R6").Value = 19144.33749835694
$ws.Range("S6").Value = 0.2176922957092124
$ws.Range("T6").Value = 0.2176922957092125
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.019504
$ws.Range("H7").Value = 54.058512
$ws.Range("I7").Value = 0.467430078412646
$ws.Range("J7").Value = 0.4674300784126461
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2342114953037475
$ws.Range("P7").Value = 0.2342114953037476
$ws.Range("Q7").Value = 1069.743447837461
$ws.Range("R7").Value = 9627.691030537151
$ws.Range("S7").Value = 0.1094774976149738
$ws.Range("T7").Value = 0.1094774976149738
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.019504
$ws.Range("H8").Value = 54.058512
$ws.Range("I8").Value = 0.467430078412646
$ws.Range("J8").Value = 0.4674300784126461
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.84506433333333
$ws.Range("N8").Value = 170.535193
$ws.Range("O8").Value = 0.2242663767030476
$ws.Range("P8").Value = 0.2242663767030477
$ws.Range("Q8").Value = 1024.319864134757
$ws.Range("R8").Value = 9218.878777212816
$ws.Range("S8").Value = 0.1048288500476256
$ws.Range("T8").Value = 0.1048288500476256
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.019504
$ws.Range("H9").Value = 54.058512
$ws.Range("I9").Value = 0.467430078412646
$ws.Range("J9").Value = 0.4674300784126461
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.21324333333333
$ws.Range("N9").Value = 57.63973
$ws.Range("O9").Value = 0.07580050295684103
$ws.Range("P9").Value = 0.07580050295684104
$ws.Range("Q9").Value = 346.2131150979734
$ws.Range("R9").Value = 3115.91803588176
$ws.Range("S9").Value = 0.03543143504083421
$ws.Range("T9").Value = 0.03543143504083422
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.752692000000001
$ws.Range("H10").Value = 26.258076
$ws.Range("I10").Value = 0.2270468436801446
$ws.Range("J10").Value = 0.2270468436801446
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 118.0470123333333
$ws.Range("N10").Value = 354.141037
$ws.Range("O10").Value = 0.4657216250363638
$ws.Range("P10").Value = 0.4657216250363638
$ws.Range("Q10").Value = 1033.229140473868
$ws.Range("R10").Value = 9299.062264264812
$ws.Range("S10").Value = 0.1057406249980942
$ws.Range("T10").Value = 0.1057406249980942
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8.752692000000001
$ws.Range("H11").Value = 26.258076
$ws.Range("I11").Value = 0.2270468436801446
$ws.Range("J11").Value = 0.2270468436801446
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2342114953037475
$ws.Range("P11").Value = 0.2342114953037476
$ws.Range("Q11").Value = 519.611134576144
$ws.Range("R11").Value = 4676.500211185296
$ws.Range("S11").Value = 0.05317698076232288
$ws.Range("T11").Value = 0.0531769807623229
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 8.752692000000001
$ws.Range("H12").Value = 26.258076
$ws.Range("I12").Value = 0.2270468436801446
$ws.Range("J12").Value = 0.2270468436801446
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.84506433333333
$ws.Range("N12").Value = 170.535193
$ws.Range("O12").Value = 0.2242663767030476
$ws.Range("P12").Value = 0.2242663767030477
$ws.Range("Q12").Value = 497.5473398298521
$ws.Range("R12").Value = 4477.926058468668
$ws.Range("S12").Value = 0.05091897297400928
$ws.Range("T12").Value = 0.05091897297400928
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 8.752692000000001
$ws.Range("H13").Value = 26.258076
$ws.Range("I13").Value = 0.2270468436801446
$ws.Range("J13").Value = 0.2270468436801446
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.21324333333333
$ws.Range("N13").Value = 57.63973
$ws.Range("O13").Value = 0.07580050295684103
$ws.Range("P13").Value = 0.07580050295684104
$ws.Range("Q13").Value = 168.16760121772
$ws.Range("R13").Value = 1513.50841095948
$ws.Range("S13").Value = 0.01721026494571822
$ws.Range("T13").Value = 0.01721026494571823
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.798175000000001
$ws.Range("H14").Value = 20.394525
$ws.Range("I14").Value = 0.1763462231431503
$ws.Range("J14").Value = 0.1763462231431503
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 118.0470123333333
$ws.Range("N14").Value = 354.141037
$ws.Range("O14").Value = 0.4657216250363638
$ws.Range("P14").Value = 0.4657216250363638
$ws.Range("Q14").Value = 802.5042480691584
$ws.Range("R14").Value = 7222.538232622425
$ws.Range("S14").Value = 0.08212824961125319
$ws.Range("T14").Value = 0.08212824961125322
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.798175000000001
$ws.Range("H15").Value = 20.394525
$ws.Range("I15").Value = 0.1763462231431503
$ws.Range("J15").Value = 0.1763462231431503
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2342114953037475
$ws.Range("P15").Value = 0.2342114953037476
$ws.Range("Q15").Value = 403.5795415624333
$ws.Range("R15").Value = 3632.2158740619
$ws.Range("S15").Value = 0.04130231261352557
$ws.Range("T15").Value = 0.04130231261352559
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.798175000000001
$ws.Range("H16").Value = 20.394525
$ws.Range("I16").Value = 0.1763462231431503
$ws.Range("J16").Value = 0.1763462231431503
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 56.84506433333333
$ws.Range("N16").Value = 170.535193
$ws.Range("O16").Value = 0.2242663767030476
$ws.Range("P16").Value = 0.2242663767030477
$ws.Range("Q16").Value = 386.4426952242584
$ws.Range("R16").Value = 3477.984257018325
$ws.Range("S16").Value = 0.03954852850958145
$ws.Range("T16").Value = 0.03954852850958146
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.798175000000001
$ws.Range("H17").Value = 20.394525
$ws.Range("I17").Value = 0.1763462231431503
$ws.Range("J17").Value = 0.1763462231431503
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.21324333333333
$ws.Range("N17").Value = 57.63973
$ws.Range("O17").Value = 0.07580050295684103
$ws.Range("P17").Value = 0.07580050295684104
$ws.Range("Q17").Value = 130.6149904975834
$ws.Range("R17").Value = 1175.53491447825
$ws.Range("S17").Value = 0.01336713240879011
$ws.Range("T17").Value = 0.01336713240879012
